# Fix validation and upload:
#  - Column A becomes plain numeric IDs (156..162) instead of text labels
#  - Column B becomes new "Model-N" text values (with some repeats),
#    replacing the old numeric-looking strings, and loses its quote-prefix
#    number-format styling (back to the workbook's default style)
#  - Selection moves to D26

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @(156, 157, 158, 159, 160, 161, 162)
$colB = @("Model-2", "Model-3", "Model-4", "Model-6", "Model-6", "Model-7", "Model-3")

for ($i = 0; $i -lt 7; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $colA[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
    # Clear the quote-prefix / number-format styling that these cells
    # previously carried, restoring the default (unstyled) cell format.
    $ws.Cells.Item($r, 2).ClearFormats()
}

$ws.Range("D26").Select()
